# add category for Binary_Search
#
# Adds a new "Binary_Search" worksheet (mirroring the existing
# LeetCode-category sheets in this workbook) after the last sheet,
# fills it with the relevant problem titles, and sorts them
# alphabetically (A-Z) the same way the other category sheets are sorted.

$wb = $excel.ActiveWorkbook

# Insert the new sheet right after the current last sheet so it lands
# at the end of the tab strip (and becomes the active tab), matching
# how Excel normally appends a worksheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Binary_Search"

# Problem titles for the new Binary_Search category.
$values = @(
    "Search Insert Position",
    "Search for a Range",
    "Sqrt(x)",
    "Search a 2D Matrix",
    "Search in Rotated Sorted Array",
    "Search in Rotated Sorted Array II"
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $newSheet.Cells.Item($i + 1, 1).Value = $values[$i]
}

# Sort column A alphabetically, same as the other category sheets.
$sortRange = $newSheet.Range("A1:A6")
$newSheet.Sort.SortFields.Clear()
$newSheet.Sort.SortFields.Add($newSheet.Range("A1")) | Out-Null
$newSheet.Sort.SetRange($sortRange)
$newSheet.Sort.Header = -4142
$newSheet.Sort.Apply()

# Match the column sizing / selection state of the other sheets.
$newSheet.Columns.Item(1).ColumnWidth = 29.14
$newSheet.Range("D8").Select() | Out-Null
